$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Insert a new "Task Three" heading paragraph (and a following
# blank paragraph) at the very top of the document, matching the formatting
# already used by the first paragraph ("Overall Solution:").
# ---------------------------------------------------------------------------
$firstParaRange = $d.Paragraphs(1).Range
$firstParaRange.InsertParagraphBefore()
$secondNewRange = $d.Paragraphs(1).Range
$secondNewRange.InsertParagraphBefore()
$d.Paragraphs(1).Range.Text = "Task Three"

# ---------------------------------------------------------------------------
# Change 2: Split "...backup attributes, we can create a new..." into three
# runs so that the word "attributes" becomes "attribute" while leaving the
# surrounding text/formatting untouched.
# ---------------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("attributes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRange.Text = "attribute"

$contextRange = $d.Content
$contextRange.Find.Execute("backup attribute, we", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$wordStart = $contextRange.Start + 7
$wordEnd = $wordStart + 9
$wordRange = $d.Range($wordStart, $wordEnd)
$wordRange.Bold = 1
$wordRange.Bold = 0

# ---------------------------------------------------------------------------
# Change 3: Remove the stray "the " run from "We can use below the tools and
# frameworks for Unit testing:" while keeping the rest of the sentence's run
# boundaries intact.
# ---------------------------------------------------------------------------
$sentenceRange = $d.Content
$sentenceRange.Find.Execute("We can use below the tools and frameworks for Unit testing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$theStart = $sentenceRange.Start + 17
$theEnd = $theStart + 4
$theRange = $d.Range($theStart, $theEnd)
$theRange.Delete()

$afterDeleteRange = $d.Content
$afterDeleteRange.Find.Execute("We can use below tools and frameworks for Unit testing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$toolsStart = $afterDeleteRange.Start + 17
$toolsEnd = $toolsStart + 20
$toolsRange = $d.Range($toolsStart, $toolsEnd)
$toolsRange.Bold = 1
$toolsRange.Bold = 0

$forStart = $toolsEnd
$forEnd = $forStart + 17
$forRange = $d.Range($forStart, $forEnd)
$forRange.Bold = 1
$forRange.Bold = 0
